# Weekly refresh: the data rows (2..36) get re-shuffled into a new row
# order (new week's rows mixed in with the existing ones, then the whole
# block was re-saved). No cell *values* actually change - only which
# row each existing record ends up on. So: snapshot all data rows as
# they are now, then write them back out in the new order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow  = 36
$lastCol      = 18   # columns A..R

# Maps destination row -> source row (both referring to the *current*,
# pre-edit layout of the sheet).
$rowMap = @{
    2  = 6
    3  = 30
    4  = 11
    5  = 36
    6  = 23
    7  = 24
    8  = 9
    9  = 12
    10 = 18
    11 = 17
    12 = 28
    13 = 21
    14 = 22
    15 = 8
    16 = 5
    17 = 29
    18 = 14
    19 = 34
    20 = 32
    21 = 20
    22 = 31
    23 = 10
    24 = 19
    25 = 7
    26 = 2
    27 = 35
    28 = 3
    29 = 4
    30 = 33
    31 = 27
    32 = 25
    33 = 13
    34 = 15
    35 = 26
    36 = 16
}

# Snapshot every data row's values first (Value2 gives back plain
# numbers/strings, e.g. the D-column date as its numeric serial,
# instead of a wrapped Variant) before writing anything back, since
# several rows are both sources and destinations in the map above.
$snapshot = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le $lastCol; $c++) {
        $rowVals += , ($ws.Cells.Item($r, $c).Value2())
    }
    $snapshot[$r] = $rowVals
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $vals = $snapshot[$srcRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $vals[$c - 1]
    }
}
